$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 1.67
$ws.Range("I2").Value = 4.33

# Row 4 updates
$ws.Range("I4").Value = 4.85
$ws.Range("J4").Value = 2.27
$ws.Range("K4").Value = 2.1
$ws.Range("N4").Value = 6.7
$ws.Range("O4").Value = 1.34
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.7
$ws.Range("W4").Value = 5.8
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 8.5
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 90
$ws.Range("AG4").Value = 900
$ws.Range("AH4").Value = 11.75
$ws.Range("AI4").Value = 28
$ws.Range("AK4").Value = 90
$ws.Range("AP4").Value = 18.5
$ws.Range("AQ4").Value = 29
$ws.Range("AR4").Value = 65
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 7.5
$ws.Range("AV4").Value = 70
$ws.Range("AX4").Value = 6.4
$ws.Range("AY4").Value = 27
$ws.Range("AZ4").Value = 32
$ws.Range("BB4").Value = 200
$ws.Range("BC4").Value = 450
